$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(312049950, Molham  Peretz: 7,8)"
$ws.Range("B1").Value = "(308073899, Anan  Kirshenbaum: 3,2)"
$ws.Range("C1").Value = "(318869187, Soaad  Leibovich: 9,2)"
$ws.Range("D1").Value = "(205898513, Asaf  Braymok: 3,6)"
$ws.Range("E1").Value = "(318428158, Tal  Asulin: 7,-7)"
$ws.Range("F1").Value = "(316028364, Sami  Castro: -3,-10)"
$ws.Range("G1").Value = "(318294931, Shalev  Afanasenko: -3,-8)"

$ws.Range("A3").Value = "cost: 389.19254997238943"
$ws.Range("A4").Value = "time: 72.83850999447787"

$wb.Save()
